$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) "Training Dashboard" sheet: update PERIOD TO EXPIRE (H) and LAST UPDATE
#    (I) for the data rows (3-21) to reflect the new "as of" date 04-Nov-2025.
#    Each PERIOD TO EXPIRE value drops by 1 day, and LAST UPDATE moves from
#    03-Nov-2025 to 04-Nov-2025.
# ---------------------------------------------------------------------------
$wsTraining = $wb.Worksheets.Item("Training Dashboard")

# Column I holds a plain text date string (not a real date) in the source
# file - force Text formatting first so COM doesn't silently convert the
# "dd-Mmm-yyyy" literal into a date serial number when we assign it.
$wsTraining.Range("I3:I21").NumberFormat = "@"

for ($row = 3; $row -le 21; $row++) {
    $periodCell = $wsTraining.Cells.Item($row, 8)      # column H
    $lastUpdateCell = $wsTraining.Cells.Item($row, 9)  # column I

    $currentPeriod = $periodCell.Value2
    if ($currentPeriod -ne $null) {
        $periodCell.Value2 = $currentPeriod - 1
    }

    $lastUpdateCell.Value2 = "04-Nov-2025"
}

# Restore the original cell style (border + no special number format) on
# column I now that the text forcing is done, pulling formats from the
# untouched column J (which keeps style index 3 throughout).
$wsTraining.Cells.Item(3, 10).Copy()
$wsTraining.Range("I3:I21").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 2) "Exam Dashboard" sheet: a new exam record ("Cs Hoist") was added as row
#    5, pushing "Climate Control Center" down to row 6 and the TOTAL AVERAGE
#    summary row down to row 7 (with an updated average).
# ---------------------------------------------------------------------------
$wsExam = $wb.Worksheets.Item("Exam Dashboard")

# Insert a new row above the current row 5 ("Climate Control Center"),
# copying the formatting of the row above (row 4) so borders/fill match.
$wsExam.Rows.Item(5).Insert()
$wsExam.Range("A4:G4").Copy()
$wsExam.Range("A5:G5").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# C (exam date) and D (marks attained) are plain text in this workbook too
# ("30-Oct-2025", "62.96%") - force Text so COM won't coerce them into a
# date serial / fractional number.
$wsExam.Range("C5:D5").NumberFormat = "@"

$wsExam.Cells.Item(5, 1).Value2 = 3
$wsExam.Cells.Item(5, 2).Value2 = "Cs Hoist"
$wsExam.Cells.Item(5, 3).Value2 = "30-Oct-2025"
$wsExam.Cells.Item(5, 4).Value2 = "62.96%"
$wsExam.Cells.Item(5, 5).Value2 = "low percentage"
$wsExam.Cells.Item(5, 6).Value2 = "This is a low mark, please retake the exam and improve your score. date is valid"

# Re-sync row 5's style to match row 4 (style index 4) now that values are
# already typed text where needed - a plain format-only paste won't disturb
# the values already entered.
$wsExam.Range("A4:G4").Copy()
$wsExam.Range("A5:G5").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# "Climate Control Center" (now row 6) keeps its data but its SN becomes 4.
$wsExam.Cells.Item(6, 1).Value2 = 4

# TOTAL AVERAGE row (now row 7): refresh the displayed average. D is a text
# percentage string here too, so guard against numeric coercion the same way.
$avgCell = $wsExam.Cells.Item(7, 4)
$avgCell.NumberFormat = "@"
$avgCell.Value2 = "65.00%"

# Restore the original style (index 3, same as the rest of the TOTAL AVERAGE
# row) by pulling format-only from its neighbor C7.
$wsExam.Cells.Item(7, 3).Copy()
$avgCell.PasteSpecial(-4122)
$excel.CutCopyMode = $false
